$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells that are purely numeric-looking strings need explicit text
# formatting first, otherwise Excel auto-converts them to numbers on
# assignment (losing the source formatting, e.g. trailing zeros).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '29.463.56'
$ws.Range("E2").Value = '  +1.04%  '
$ws.Range("D3").Value = '1.839.73'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("D4").Value = '0.9984'
$ws.Range("E4").Value = '  -0.88%  '
$ws.Range("D5").Value = '243.56'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").Value = '0.6271'
$ws.Range("E6").Value = '  +1.62%  '
$ws.Range("D7").Value = '0.9996'
$ws.Range("E7").Value = '  -0.96%  '
$ws.Range("D8").Value = '0.07422'
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '0.2942'
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("D10").Value = '23.64'
$ws.Range("E10").Value = '  +3.24%  '
$ws.Range("D11").Value = '0.07656'
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("D12").Value = '1.837.63'
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").Value = '5.017'
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").Value = '0.6770'
$ws.Range("E14").Value = '  +1.21%  '
$ws.Range("D15").Value = '83.56'
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").Value = '0.000009283'
$ws.Range("E16").Value = '  +2.31%  '
$ws.Range("D17").Value = '5.919'
$ws.Range("E17").Value = '  +0.84%  '
$ws.Range("D18").Value = '29.422.90'
$ws.Range("E18").Value = '  +0.90%  '
$ws.Range("D19").Value = '2.081.94'
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").Value = '237.47'
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("D21").Value = '12.57'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '0.9992'
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("D23").Value = '7.341'
$ws.Range("E23").Value = '  +2.74%  '
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.92%  '
$ws.Range("D25").Value = '159.08'
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("D26").Value = '0.1412'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '8.504'
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").Value = '17.75'
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '1.496'
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '0.06026'
$ws.Range("E30").Value = '  +8.56%  '
$ws.Range("D31").Value = '1.238'
$ws.Range("E31").Value = '  +1.88%  '
$ws.Range("D32").Value = '4.098'
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("D33").Value = '4.111'
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("D34").Value = '1.878'
$ws.Range("E34").Value = '  +2.15%  '
$ws.Range("D35").Value = '1.142'
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("D36").Value = '0.7268'
$ws.Range("E36").Value = '  -1.87%  '
$ws.Range("E37").Value = '  -1.47%  '
$ws.Range("D38").Value = '2.880'
$ws.Range("E38").Value = '  +2.07%  '
$ws.Range("D39").Value = '1.218.26'
$ws.Range("E39").Value = '  +1.21%  '
$ws.Range("D40").Value = '0.01761'
$ws.Range("E40").Value = '  -0.74%  '
$ws.Range("D41").Value = '6.265'
$ws.Range("E41").Value = '  -2.29%  '
$ws.Range("D42").Value = '0.9107'
$ws.Range("E42").Value = '  +1.22%  '
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("D44").Value = '1.997.99'
$ws.Range("E44").Value = '  +0.87%  '
$ws.Range("D45").Value = '102.01'
$ws.Range("E45").Value = '  +0.57%  '
$ws.Range("D46").Value = '65.50'
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("E48").Value = '  -1.07%  '
$ws.Range("D49").Value = '9.230'
$ws.Range("E49").Value = '  +1.28%  '
$ws.Range("D50").Value = '0.4063'
$ws.Range("E50").Value = '  +0.97%  '
$ws.Range("E51").Value = '  +3.19%  '
